$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Kyrgyz title in A1 (row 1) ---
$ws.Range("A1").Value = "4.2.2.1а Балдарды мектепке чейин билим берүү менен камтуу"

# --- 2. Update "urban settlements" -> "urban" (English column C, urban rows) ---
$urbanRows = @(6, 9, 12, 15, 18, 21, 24, 27)
foreach ($r in $urbanRows) {
    $ws.Cells.Item($r, 3).Value = "urban"
}

# --- 3. Update "countryside" -> "rural" (English column C, rural rows) ---
$ruralRows = @(7, 10, 13, 16, 19, 22, 25, 28)
foreach ($r in $ruralRows) {
    $ws.Cells.Item($r, 3).Value = "rural"
}

# --- 4. Add new column N with 2023 data ---
# Header row 4: year 2023, same style as M4
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# Data rows 5-29: copy number format from the corresponding M-column cell, then
# apply the percentage-style numeric format (0.0) used by the rest of the table,
# and finally write the 2023 value.
$values = @{
    5  = 28.34784779265912
    6  = 39.999446500300472
    7  = 23.198557483143556
    8  = 27.597876990321573
    9  = 47.175678010018999
    10 = 22.17579894112394
    11 = 24.100104034215697
    12 = 38.296287676015361
    13 = 19.410249509822766
    14 = 30.400174646089773
    15 = 44.562134629854725
    16 = 24.612036336109007
    17 = 39.266683582846994
    18 = 54.818496110630946
    19 = 36.591078066914498
    20 = 23.890520476423561
    21 = 16.93085228577992
    22 = 24.386979772654026
    23 = 28.919699950811605
    24 = 37.932834522359492
    25 = 26.985549456704376
    26 = 27.190143693828379
    27 = 54.006768771869439
    28 = 22.334624692306893
    29 = 36.01461582008131
}

foreach ($r in 5..29) {
    $ws.Range("M$r").Copy()
    $ws.Range("N$r").PasteSpecial(-4122)
    $ws.Range("N$r").NumberFormat = "0.0"
    $ws.Range("N$r").Value = $values[$r]
}

# Row 30 (bottom total row) already carries the percentage format on M30, so a
# straight format copy (no NumberFormat override needed) matches the rest.
$ws.Range("M30").Copy()
$ws.Range("N30").PasteSpecial(-4122)
$ws.Range("N30").Value = 42.081208505725009

# Reset selection back to the default cell.
$ws.Range("A1").Select()
